$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44558
$ws.Range("M2").Value = 20

# Row 3
$ws.Range("D3").Value = 44558
$ws.Range("M3").Value = 25

# Row 5
$ws.Range("D5").Value = 44572
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("Q5").Value = "`$/bandeja 6 kilos"
$ws.Range("R5").Value = "Región Metropolitana"
$ws.Range("S5").Value = 3333
$ws.Range("T5").Value = 6

# Row 6
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 20
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 2143

# Row 7
$ws.Range("D7").Value = 44189
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 13000
$ws.Range("O7").Value = 13000
$ws.Range("P7").Value = 13000
$ws.Range("R7").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S7").Value = 1857

# Row 8
$ws.Range("D8").Value = 44553
$ws.Range("M8").Value = 200

# Row 9
$ws.Range("D9").Value = 44553
$ws.Range("M9").Value = 150

# Row 10
$ws.Range("D10").Value = 44204
$ws.Range("M10").Value = 110
$ws.Range("N10").Value = 7000
$ws.Range("O10").Value = 7500
$ws.Range("P10").Value = 7318
$ws.Range("Q10").Value = "`$/bandeja 7 kilos"
$ws.Range("R10").Value = "Provincia de San Felipe de Aconcagua"
$ws.Range("S10").Value = 1045
$ws.Range("T10").Value = 7

# Row 11
$ws.Range("D11").Value = 44187
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 45
$ws.Range("N11").Value = 14000
$ws.Range("O11").Value = 14000
$ws.Range("P11").Value = 14000
$ws.Range("S11").Value = 2000

# Row 12
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 50
$ws.Range("N12").Value = 12000
$ws.Range("O12").Value = 12000
$ws.Range("P12").Value = 12000
$ws.Range("S12").Value = 1714

# Row 13
$ws.Range("D13").Value = 44550
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 24000
$ws.Range("O13").Value = 24000
$ws.Range("P13").Value = 24000
$ws.Range("R13").Value = "Región Metropolitana"
$ws.Range("S13").Value = 3429
